$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Row 10: put "DepartureAirport'," in B10 (style s2 / yellow) and restyle C10 to s2 (value unchanged).
$ws.Range("B10").Value2 = $ws.Range("C10").Value2
$ws.Range("B10").Interior.Color = $ws.Range("A7").Interior.Color
$ws.Range("C10").Interior.Color = $ws.Range("A7").Interior.Color

# 2. Delete row 15 entirely (it held the now-duplicate "DepartureAirport'," entry),
#    shifting everything below it up by one row.
$ws.Rows("15").Delete()

# 3. Row 20 (was row 21 before the shift): "PassengerBirthDate'," -- recolor to yellow (s2).
$ws.Range("A20:B20").Interior.Color = $ws.Range("A7").Interior.Color

# 4. Row 22 (was row 23 before the shift): "PassengerDocument',"
#    recolor A22/B22 from green to yellow, and clear C22 (was green "PassengerDocument'" no-comma).
$ws.Range("A22:B22").Interior.Color = $ws.Range("A7").Interior.Color
$ws.Range("C22").ClearContents()
$ws.Range("C22").Interior.Pattern = -4142

$ws.Range("A28").Select()
